$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture original comment texts for cells that will shift position (G1,H1,I1 -> J1,K1,L1) ---
$textG1 = $ws.Range("G1").Comment.Text()
$textH1 = $ws.Range("H1").Comment.Text()
$textI1 = $ws.Range("I1").Comment.Text()

# Remove the comments that will be relocated; they'll be re-added at their new cells below.
$ws.Range("G1").Comment.Delete()
$ws.Range("H1").Comment.Delete()
$ws.Range("I1").Comment.Delete()

# --- Insert three new columns before the (old) G column, shifting G,H,I -> J,K,L ---
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").EntireColumn.Insert()

# The insert drags blank styled cells into the new columns for every existing data row;
# clear them back out so rows 2-9 don't carry stray empty G:I cells.
$ws.Range("G2:I9").Clear()

# --- Fix up the Verified-column data validation (column is now K, was H) ---
$ws.Columns("K").Validation.Delete()
$ws.Range("K2:K183").Validation.Add(3, 1, 1, '"Yes,No"')

# --- Match the existing header/body styling for the 3 new fund-currency columns ---
$ws.Range("F1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# --- Re-label the Folio-currency columns (D:F) and label the new Fund-currency columns (G:I) ---
$ws.Range("D1").Value = "Call Amount (Inclusive of Capital Fees, Folio Currency)*"
$ws.Range("E1").Value = "Capital Fees (Folio Currency)"
$ws.Range("F1").Value = "Other Fees (Folio Currency)"
$ws.Range("G1").Value = "Call Amount (Inclusive of Capital Fees, Fund Currency)*"
$ws.Range("H1").Value = "Capital Fees (Fund Currency)"
$ws.Range("I1").Value = "Other Fees (Fund Currency)"

# --- New fund-currency amount for row 2 ---
$ws.Range("G2").Value = 1200000

# --- Re-create the relocated comments at their new homes ---
$ws.Range("J1").AddComment($textG1) | Out-Null
$ws.Range("K1").AddComment($textH1) | Out-Null
$ws.Range("L1").AddComment($textI1) | Out-Null

# --- Brand-new comments for the new fund-currency columns (same guidance as folio-currency cols) ---
$ws.Range("G1").AddComment("Author:`nAmount being called.  Include any fees amount which is part of the commiment amount as well") | Out-Null
$ws.Range("H1").AddComment("Author:`nApplicable in case there are Fees / expenses that are part of the commitment amount which need to be tracked as such") | Out-Null
$ws.Range("I1").AddComment("Author:`nApplicable in case there are fees / expenses over and above the Commitment amount being called ") | Out-Null

# --- New data row (Investor 6 / row 9), formatted the same as row 8 ---
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
$ws.Range("K8:L8").Copy()
$ws.Range("K9:L9").PasteSpecial(-4122)

$ws.Range("A9").Value = "Investor 6"
$ws.Range("B9").Value = "SAAS Fund"
$ws.Range("C9").Value = "Call 1"
$ws.Range("D9").Value = 13000
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = 5
$ws.Range("K9").Value = "Yes"
$ws.Range("L9").Value = 13

# --- Selection as left by the edit ---
$ws.Range("L9").Select()

Write-Output "done"
